$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old "Unnamed: 0.x" columns (B..U) entirely, shifting the
# surviving "our_*" columns (previously V..Y) left into B..E.
$ws.Range("B1:U1").EntireColumn.Delete()

# Update the header row (row 1) with the new column names.
$ws.Range("B1").Value = "our_identified"
$ws.Range("C1").Value = "our_Overlap_merlin"
$ws.Range("D1").Value = "ourbest_param"
$ws.Range("E1").Value = "ourtime_taken"

# Update the data row (row 2) with the new values.
$ws.Range("B2").Value = "[]"
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = "{'cluster': 26, 'training': 386, 'window': 362, 'threshold': 1.5}"
$ws.Range("E2").Value = 20.03201633800199
